# The commit changes the value of cell C11 on the "Rules" sheet
# from 22 to 2211 (numeric), leaving its style/format untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = 2211
